$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header: column A title changes from "Gen" to "MaxFES"
$ws.Range("A1").Value = "MaxFES"

# Column A (MaxFES fraction) values for rows 3-14
$ws.Range("A3").Value = 0.001
$ws.Range("A4").Value = 0.01
$ws.Range("A5").Value = 0.1
$ws.Range("A6").Value = 0.2
$ws.Range("A7").Value = 0.3
$ws.Range("A8").Value = 0.4
$ws.Range("A9").Value = 0.5
$ws.Range("A10").Value = 0.6
$ws.Range("A11").Value = 0.7
$ws.Range("A12").Value = 0.8
$ws.Range("A13").Value = 0.9
$ws.Range("A14").Value = 1

# New AZ column values (previously "Run 50" data col, now becomes "Mean" col after BA removal)
$ws.Range("AZ2").Value = 140.39512052
$ws.Range("AZ3").Value = 123.28081204
$ws.Range("AZ4").Value = 70.41294789
$ws.Range("AZ5").Value = 22.28285643
$ws.Range("AZ6").Value = 14.51282172
$ws.Range("AZ7").Value = 10.98339361
$ws.Range("AZ8").Value = 8.88507302
$ws.Range("AZ9").Value = 7.68544035
$ws.Range("AZ10").Value = 6.62948054
$ws.Range("AZ11").Value = 5.73737988
$ws.Range("AZ12").Value = 4.99504479
$ws.Range("AZ13").Value = 4.42704347
$ws.Range("AZ14").Value = 4.10683071

# AZ1 header becomes "Mean" (it was "Run 50" at col index 51; Mean string removed from BA col)
$ws.Range("AZ1").Value = "Mean"

# Delete the BA column entirely (previously held the Mean data, now removed)
$ws.Columns("BA").Delete()
